$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2994
$ws.Range("I125").Value = 4334.1665
$ws.Range("J125").Value = 1845.2858
$ws.Range("K125").Value = 39007.4985
$ws.Range("L125").Value = 16607.5722
$ws.Range("M125").Value = -36547.4985
$ws.Range("N125").Value = -21527.5722

$ws.Range("H131").Value = 1992.6
$ws.Range("I131").Value = 1992.6
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 5977.799999999999
$ws.Range("L131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -937.7999999999993

$ws.Range("H138").Value = 4658.3955
$ws.Range("I138").Value = 1886.091
$ws.Range("J138").Value = 5611.375
$ws.Range("K138").Value = 5658.272999999999
$ws.Range("L138").Value = 16834.125
$ws.Range("M138").Value = -518.2729999999992
$ws.Range("N138").Value = -27114.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 44949.75
$ws.Range("J24").Value = 44949.75
$ws.Range("L24").Value = 44949.75
$ws.Range("N24").Value = -45697.75

$ws.Range("H32").Value = 1951159.4
$ws.Range("I32").Value = 916.25
$ws.Range("J32").Value = 12352456
$ws.Range("K32").Value = 916.25
$ws.Range("L32").Value = 12352456
$ws.Range("M32").Value = -629.25
$ws.Range("N32").Value = -12353030

$ws.Range("H100").Value = 44949.75
$ws.Range("J100").Value = 44949.75
$ws.Range("L100").Value = 44949.75
$ws.Range("N100").Value = -47113.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H125").Value = 100001
$ws.Range("J125").Value = 100001
$ws.Range("L125").Value = 100001
$ws.Range("N125").Value = -109841

$ws.Range("H141").Value = 220000
$ws.Range("J141").Value = 220000
$ws.Range("L141").Value = 220000
$ws.Range("N141").Value = -230360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5143.143
$ws.Range("I31").Value = 2093.182
$ws.Range("J31").Value = 7116.647
$ws.Range("K31").Value = 2093.182
$ws.Range("L31").Value = 7116.647
$ws.Range("M31").Value = -1798.182
$ws.Range("N31").Value = -7706.647

$ws.Range("H34").Value = 5143.143
$ws.Range("I34").Value = 2093.182
$ws.Range("J34").Value = 7116.647
$ws.Range("K34").Value = 2093.182
$ws.Range("L34").Value = 7116.647
$ws.Range("M34").Value = -1891.182
$ws.Range("N34").Value = -7520.647

$ws.Range("H58").Value = 25647028
$ws.Range("I58").Value = 37041350
$ws.Range("J58").Value = 9797.666999999999
$ws.Range("K58").Value = 37041350
$ws.Range("L58").Value = 9797.666999999999
$ws.Range("M58").Value = -37041147
$ws.Range("N58").Value = -10203.667

$ws.Range("H136").Value = 25647028
$ws.Range("I136").Value = 37041350
$ws.Range("J136").Value = 9797.666999999999
$ws.Range("K136").Value = 111124050
$ws.Range("L136").Value = 29393.001
$ws.Range("M136").Value = -111121500
$ws.Range("N136").Value = -34493.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 18723.445
$ws.Range("J62").Value = 19812.375
$ws.Range("L62").Value = 59437.125
$ws.Range("N62").Value = -60809.125

$ws.Range("H65").Value = 18723.445
$ws.Range("J65").Value = 19812.375
$ws.Range("L65").Value = 178311.375
$ws.Range("N65").Value = -185175.375

$ws.Range("H70").Value = 4900
$ws.Range("I70").Value = 4900
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 14700
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -14385

$ws.Range("H73").Value = 4900
$ws.Range("I73").Value = 4900
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 14700
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -13608

$ws.Range("H74").Value = 13619
$ws.Range("I74").Value = 4738.5
$ws.Range("J74").Value = 22499.5
$ws.Range("K74").Value = 14215.5
$ws.Range("L74").Value = 67498.5
$ws.Range("M74").Value = -13154.5
$ws.Range("N74").Value = -69620.5

$ws.Range("H77").Value = 13619
$ws.Range("I77").Value = 4738.5
$ws.Range("J77").Value = 22499.5
$ws.Range("K77").Value = 42646.5
$ws.Range("L77").Value = 202495.5
$ws.Range("M77").Value = -37342.5
$ws.Range("N77").Value = -213103.5

$ws.Range("H82").Value = 18002.4
$ws.Range("I82").Value = 7006.5
$ws.Range("J82").Value = 25333
$ws.Range("K82").Value = 21019.5
$ws.Range("L82").Value = 75999
$ws.Range("M82").Value = -20613.5
$ws.Range("N82").Value = -76811

$ws.Range("H85").Value = 18002.4
$ws.Range("I85").Value = 7006.5
$ws.Range("J85").Value = 25333
$ws.Range("K85").Value = 21019.5
$ws.Range("L85").Value = 75999
$ws.Range("M85").Value = -19615.5
$ws.Range("N85").Value = -78807

$ws.Range("H87").Value = 23279.3
$ws.Range("I87").Value = 18978.6
$ws.Range("J87").Value = 27580
$ws.Range("K87").Value = 56935.8
$ws.Range("L87").Value = 82740
$ws.Range("M87").Value = -55687.8
$ws.Range("N87").Value = -85236

$ws.Range("H90").Value = 23279.3
$ws.Range("I90").Value = 18978.6
$ws.Range("J90").Value = 27580
$ws.Range("K90").Value = 170807.4
$ws.Range("L90").Value = 248220
$ws.Range("M90").Value = -164567.4
$ws.Range("N90").Value = -260700

$ws.Range("H93").Value = 6200
$ws.Range("J93").Value = 1500
$ws.Range("L93").Value = 4500
$ws.Range("N93").Value = -8244

$ws.Range("H95").Value = 7242
$ws.Range("J95").Value = 9994
$ws.Range("L95").Value = 29982
$ws.Range("N95").Value = -34100

$ws.Range("H96").Value = 18399
$ws.Range("J96").Value = 18399
$ws.Range("L96").Value = 55197
$ws.Range("N96").Value = -59315

$ws.Range("H99").Value = 3524
$ws.Range("I99").Value = 3524
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 10572
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -8326

$ws.Range("H100").Value = 22750
$ws.Range("J100").Value = 24500
$ws.Range("L100").Value = 73500
$ws.Range("N100").Value = -75122

$ws.Range("H105").Value = 29999.375
$ws.Range("J105").Value = 29999.375
$ws.Range("L105").Value = 89998.125
$ws.Range("N105").Value = -95240.125

$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 90000
$ws.Range("N106").Value = -91892

$ws.Range("H108").Value = 6176.923
$ws.Range("I108").Value = 3200
$ws.Range("J108").Value = 7500
$ws.Range("K108").Value = 9600
$ws.Range("L108").Value = 22500
$ws.Range("M108").Value = -6720
$ws.Range("N108").Value = -28260

$ws.Range("H110").Value = 21220.445
$ws.Range("I110").Value = 11621
$ws.Range("J110").Value = 28900
$ws.Range("K110").Value = 34863
$ws.Range("L110").Value = 86700
$ws.Range("M110").Value = -30773
$ws.Range("N110").Value = -94880

$ws.Range("H111").Value = 24999
$ws.Range("I111").Value = 24999
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 74997
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -71930

$ws.Range("H112").Value = 24426.908
$ws.Range("I112").Value = 27589.5
$ws.Range("J112").Value = 15993.333
$ws.Range("K112").Value = 82768.5
$ws.Range("L112").Value = 47979.999
$ws.Range("M112").Value = -81660.5
$ws.Range("N112").Value = -50195.999

$ws.Range("H118").Value = 1447.75
$ws.Range("I118").Value = 1447.75
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 4343.25
$ws.Range("L118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = -3100.25

$ws.Range("H119").Value = 3311.4
$ws.Range("I119").Value = 3311.4
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 9934.200000000001
$ws.Range("L119").Value = 0
$ws.Range("M119").ClearContents()
$ws.Range("N119").Value = -5096.200000000001

$ws.Range("H120").Value = 30678.625
$ws.Range("I120").Value = 11357.25
$ws.Range("J120").Value = 50000
$ws.Range("K120").Value = 34071.75
$ws.Range("L120").Value = 150000
$ws.Range("M120").Value = -29233.75
$ws.Range("N120").Value = -159676

$ws.Range("H123").Value = 4344
$ws.Range("J123").Value = 5001
$ws.Range("L123").Value = 15003
$ws.Range("N123").Value = -19903

$ws.Range("H124").Value = 15318
$ws.Range("I124").Value = 6647.5
$ws.Range("J124").Value = 50000
$ws.Range("K124").Value = 19942.5
$ws.Range("L124").Value = 150000
$ws.Range("M124").Value = -15032.5
$ws.Range("N124").Value = -159820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 24874.75
$ws.Range("J98").Value = 24874.75
$ws.Range("L98").Value = 24874.75
$ws.Range("N98").Value = -30864.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 34999.5
$ws.Range("J127").Value = 34999.5
$ws.Range("L127").Value = 34999.5
$ws.Range("N127").Value = -44919.5

$ws.Range("H132").Value = 3417.42
$ws.Range("I132").Value = 2330.2559
$ws.Range("K132").Value = 6990.7677
$ws.Range("M132").Value = -4460.7677

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 52083.332
$ws.Range("I56").Value = 70000
$ws.Range("J56").Value = 43125
$ws.Range("K56").Value = 70000
$ws.Range("L56").Value = 43125
$ws.Range("M56").Value = -69286
$ws.Range("N56").Value = -44553

$ws.Range("H82").Value = 98332.664
$ws.Range("J82").Value = 98332.664
$ws.Range("L82").Value = 98332.664
$ws.Range("N82").Value = -99098.664

$ws.Range("H85").Value = 98332.664
$ws.Range("J85").Value = 98332.664
$ws.Range("L85").Value = 98332.664
$ws.Range("N85").Value = -100984.664

$ws.Range("H88").Value = 45000
$ws.Range("J88").Value = 45000
$ws.Range("L88").Value = 45000
$ws.Range("N88").Value = -45812

$ws.Range("H91").Value = 45000
$ws.Range("J91").Value = 45000
$ws.Range("L91").Value = 45000
$ws.Range("N91").Value = -47808

$ws.Range("H136").Value = 6670733.5
$ws.Range("I136").Value = 8197444.5
$ws.Range("J136").Value = 18635
$ws.Range("K136").Value = 24592333.5
$ws.Range("L136").Value = 55905
$ws.Range("M136").Value = -24589783.5
$ws.Range("N136").Value = -61005
